# Hortaliza, Terminal La Palmera de La Serena - Betarraga
# Insert one new weekly pair of records (Primera/Segunda) at row 140,
# pushing the existing rows 140:213 down to 142:215.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the existing date-block (row 140),
# shifting the rest of the table (old 140:213) down to 142:215.
$ws.Rows("140:141").Insert()

# Row 140 - "Primera" quality record for the new date.
$ws.Range("A140").Value = 8
$ws.Range("B140").Value = "Terminal La Palmera de La Serena"
$ws.Range("C140").Value = "Coquimbo"
$ws.Range("D140").Value = 44553
$ws.Range("E140").Value = 4
$ws.Range("F140").Value = 100114014
$ws.Range("G140").Value = "Betarraga"
$ws.Range("H140").Value = "Sin especificar"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 2960
$ws.Range("K140").Value = 450
$ws.Range("L140").Value = 500
$ws.Range("M140").Value = 475
$ws.Range("N140").Value = "$/paquete 3 unidades"
$ws.Range("O140").Value = "Provincia del Elquí"
$ws.Range("P140").Value = 158
$ws.Range("Q140").Value = 3
$ws.Range("R140").Value = "Hortaliza"

# Row 141 - "Segunda" quality record for the same new date.
$ws.Range("A141").Value = 8
$ws.Range("B141").Value = "Terminal La Palmera de La Serena"
$ws.Range("C141").Value = "Coquimbo"
$ws.Range("D141").Value = 44553
$ws.Range("E141").Value = 4
$ws.Range("F141").Value = 100114014
$ws.Range("G141").Value = "Betarraga"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Segunda"
$ws.Range("J141").Value = 1400
$ws.Range("K141").Value = 350
$ws.Range("L141").Value = 400
$ws.Range("M141").Value = 375
$ws.Range("N141").Value = "$/paquete 3 unidades"
$ws.Range("O141").Value = "Provincia del Elquí"
$ws.Range("P141").Value = 125
$ws.Range("Q141").Value = 3
$ws.Range("R141").Value = "Hortaliza"
